$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 0.9489610112883385, 0.1758320240729323, -0.7249361535492151, 0.6700277771922065, 1.570796383038867, -0.6218353262550272),
    @(0.206526356614631, 0.9506995464141992, 0.1758735022814049, -0.7248205873556631, 0.6701018665026294, 1.570796382838802, -0.6200967910929992),
    @(0.4130527132292621, 0.961541846730508, 0.1761321793631349, -0.7240998637083433, 0.6705639213336673, 1.570796381591107, -0.6092544905511341),
    @(0.6195790698438931, 0.9871833712480381, 0.1767439383570568, -0.7223953865647509, 0.6716566590307295, 1.570796378640369, -0.5836129655001743),
    @(0.8261054264585241, 1.03009065733366, 0.1777676262485357, -0.7195431970088403, 0.67348519340482, 1.570796373702746, -0.5407056785219373),
    @(1.032631783073155, 1.090039474317657, 0.1791978933116294, -0.7155582003966157, 0.6760399686549388, 1.570796366804042, -0.4807568602908032),
    @(1.239158139687786, 1.164652977101047, 0.1809780324513496, -0.7105983935017208, 0.6792196932904817, 1.570796358217776, -0.4061433559552025),
    @(1.445684496302417, 1.249939859762893, 0.1830128185459232, -0.7049290916610294, 0.6828542740536421, 1.570796348403254, -0.3208564715191031),
    @(1.652210852917048, 1.340832509167626, 0.1851813477890544, -0.6988871559202345, 0.6867277498418117, 1.57079633794364, -0.2299638202234981),
    @(1.858737209531679, 1.431725158572359, 0.1873498770321855, -0.6928452201794398, 0.6906012256299814, 1.570796327484026, -0.1390711689278932),
    @(2.06526356614631, 1.517012041234205, 0.1893846631267591, -0.6871759183387484, 0.6942358063931418, 1.570796317669504, -0.05378428449179418),
    @(2.271789922760941, 1.591625544017594, 0.1911648022664793, -0.6822161114438535, 0.6974155310286847, 1.570796309083238, 0.0208292198438067),
    @(2.478316279375572, 1.651574361001592, 0.192595069329573, -0.6782311148316289, 0.6999703062788034, 1.570796302184534, 0.08077803807494088),
    @(2.684842635990203, 1.694481647087213, 0.193618757221052, -0.6753789252757183, 0.701798840652894, 1.570796297246911, 0.1236853250531775),
    @(2.891368992604835, 1.720123171604743, 0.1942305162149738, -0.6736744481321258, 0.7028915783499562, 1.570796294296173, 0.1493268501041367),
    @(3.097895349219466, 1.730965471921052, 0.1944891932967038, -0.6729537244848061, 0.7033536331809941, 1.570796293048478, 0.1601691506460023),
    @(3.304421705834097, 1.732704007046913, 0.1945306715051764, -0.672838158291254, 0.7034277224914169, 1.570796292848413, 0.1619076858080303)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}